$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C5").Value = 45212
